$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.445647641019636
$ws.Range("C2").Value = 9.983522426115931
$ws.Range("D2").Value = 186123.597850132
$ws.Range("E2").Value = 2797.565817734744
$ws.Range("G2").Value = 188932.5928379339
